# Regenerate save_data "K" column (col G) values for mayza_tim.xlsx
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (column G, the 7th column)
$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 1
    39 = 0
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 1
    48 = 0
    49 = 0
    50 = 0
    51 = 1
    52 = 0
    53 = 1
    54 = 2
    55 = 0
    56 = 0
    57 = 0
    58 = 1
    59 = 1
    60 = 1
    61 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
